$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix the spelling of "delaugney" -> "Delaunay" in the Lane Generation
#    section ("Use delaugney triangulation to find edges.").
#    We use Find/Replace scoped to just the misspelled word so the
#    surrounding "Use " / " triangulation to find edges." runs are left
#    completely untouched.
# ---------------------------------------------------------------------------
$rFix = $d.Content
$rFix.Find.Execute("delaugney", $true, $false, $false, $false, $false, $true, 1, $false, "Delaunay", 2)

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark. It currently sits after "Zoom levels" --
#    it needs to end up after "Just increased sizes and power for now"
#    instead (reflecting where the author's cursor was when the document
#    was last saved). Adding a bookmark with a name that already exists
#    moves it, so we just add "_GoBack" at the new location; Word drops the
#    old one automatically.
#
#    A zero-length (collapsed) range right at the end of a paragraph isn't
#    accepted when adding a bookmark, so we briefly insert a placeholder
#    character after the target text, anchor the bookmark to that
#    character, and then delete the placeholder again -- the bookmark
#    stays behind at the correct (now collapsed) position.
# ---------------------------------------------------------------------------
$rGoBack = $d.Content
$rGoBack.Find.Execute("Just increased sizes and power for now")
$rGoBack.Collapse(0)
$rGoBack.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $rGoBack)
$rGoBack.Text = ""

Write-Output "done"
